$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "2" = "  +0.10%  "
    "3" = "  +1.71%  "
    "4" = "  -0.20%  "
    "5" = "  +1.37%  "
    "6" = "  +0.83%  "
    "7" = "  -0.19%  "
    "8" = "  +0.00%  "
    "9" = "  -4.97%  "
    "10" = "  -0.57%  "
    "11" = "  +3.53%  "
    "12" = "  +1.44%  "
    "13" = "  -3.26%  "
    "14" = "  +1.32%  "
    "15" = "  +3.63%  "
    "16" = "  +0.15%  "
    "17" = "  -0.74%  "
    "18" = "  +1.45%  "
    "19" = "  +4.94%  "
    "20" = "  +3.01%  "
    "21" = "  +0.99%  "
    "22" = "  +0.22%  "
    "23" = "  +1.80%  "
    "24" = "  +0.98%  "
    "25" = "  +1.34%  "
    "26" = "  +0.05%  "
    "27" = "  -2.06%  "
    "28" = "  -0.86%  "
    "29" = "  -0.15%  "
    "30" = "  -0.61%  "
    "31" = "  -0.12%  "
    "32" = "  +1.17%  "
    "33" = "  +1.43%  "
    "34" = "  -1.06%  "
    "35" = "  -0.07%  "
    "36" = "  +0.78%  "
    "37" = "  -0.05%  "
    "38" = "  +2.57%  "
    "39" = "  +0.99%  "
    "40" = "  -1.10%  "
    "41" = "  +1.51%  "
    "42" = "  +2.93%  "
    "43" = "  -0.20%  "
    "44" = "  -0.51%  "
    "45" = "  -0.81%  "
    "46" = "  +1.24%  "
    "47" = "  -3.21%  "
    "48" = "  +2.15%  "
    "49" = "  +4.31%  "
    "50" = "  +9.51%  "
    "51" = "  +0.48%  "
}

foreach ($row in $updates.Keys) {
    $ws.Range("E$row").Value = $updates[$row]
}
